$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the justification text in B10 and B11: "DCA" -> "ACD"
$ws.Range("B10").Value = '$\angle ACD >  \angle ABC$'
$ws.Range("B11").Value = '$\angle ACD >  \angle ABC$'

# Move the active selection to B11
$ws.Range("B11").Select()
